$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-05-21 Wednesday" "2025-05-22 Thursday"

Replace-Text "615×5=" "524×8="
Replace-Text "128×2=" "572×3="
Replace-Text "682×4=" "263×5="
Replace-Text "131×6=" "846×6="
Replace-Text "774×3=" "363×8="
Replace-Text "700×5=" "734×2="
Replace-Text "670×4=" "992×9="
Replace-Text "819×7=" "989×9="
Replace-Text "900×8=" "747×4="
Replace-Text "289×9=" "432×2="
Replace-Text "720×8=" "895×7="
Replace-Text "362×8=" "985×6="
Replace-Text "698×6=" "389×7="
Replace-Text "970×2=" "428×9="
Replace-Text "923×5=" "386×5="
Replace-Text "594×5=" "751×3="
Replace-Text "621×2=" "971×9="
Replace-Text "914×2=" "657×7="
Replace-Text "631×3=" "856×9="
Replace-Text "490×8=" "389×3="
Replace-Text "959×8=" "444×9="
Replace-Text "404×7=" "261×2="
Replace-Text "306×8=" "475×2="
Replace-Text "764×5=" "873×8="
Replace-Text "179×7=" "146×5="
